$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (swap/results updates) ---
# Row 233
$ws.Range('A233').Value = 231
$ws.Range('B233').Value = 7559468
$ws.Range('C233').Value = 'Uruguay Primera División'
$ws.Range('D233').Value = 'Uruguay Clausura'
$ws.Range('E233').Value = 45266.70833333334
$ws.Range('F233').Value = 'Liverpool Montevideo'
$ws.Range('G233').Value = 'CA River Plate'
$ws.Range('H233').Value = 2
$ws.Range('I233').Value = 1
$ws.Range('J233').Value = 'H'
$ws.Range('K233').Value = 1.7
$ws.Range('L233').Value = 3
$ws.Range('M233').Value = 5.75
$ws.Range('N233').Value = 1.833
$ws.Range('O233').Value = 3.2
$ws.Range('P233').Value = 4.5
$ws.Range('Q233').Value = -0.5
$ws.Range('R233').Value = 1.925
$ws.Range('S233').Value = 1.925
$ws.Range('T233').Value = 2.25
$ws.Range('U233').Value = 2.025
$ws.Range('V233').Value = 1.825
$ws.Range('W233').Value = 0.833
$ws.Range('X233').Value = -1
$ws.Range('Y233').Value = -1
$ws.Range('Z233').Value = 0.925
$ws.Range('AA233').Value = -1
$ws.Range('AB233').Value = 1.025
$ws.Range('AC233').Value = -1

# Row 234
$ws.Range('A234').Value = 232
$ws.Range('B234').Value = 7559469
$ws.Range('C234').Value = 'Uruguay Primera División'
$ws.Range('D234').Value = 'Uruguay Clausura'
$ws.Range('E234').Value = 45266.70833333334
$ws.Range('F234').Value = 'Montevideo Wanderers'
$ws.Range('G234').Value = 'Penarol'
$ws.Range('H234').Value = 0
$ws.Range('I234').Value = 0
$ws.Range('J234').Value = 'D'
$ws.Range('K234').Value = 4.75
$ws.Range('L234').Value = 3.4
$ws.Range('M234').Value = 1.7
$ws.Range('N234').Value = 2.7
$ws.Range('O234').Value = 3.2
$ws.Range('P234').Value = 2.45
$ws.Range('Q234').Value = 0
$ws.Range('R234').Value = 2.05
$ws.Range('S234').Value = 1.8
$ws.Range('T234').Value = 2.5
$ws.Range('U234').Value = 1.975
$ws.Range('V234').Value = 1.875
$ws.Range('W234').Value = -1
$ws.Range('X234').Value = 2.2
$ws.Range('Y234').Value = -1
$ws.Range('Z234').Value = 0
$ws.Range('AA234').Value = -0
$ws.Range('AB234').Value = -1
$ws.Range('AC234').Value = 0.875

# Row 236
$ws.Range('A236').Value = 234
$ws.Range('B236').Value = 7013409
$ws.Range('C236').Value = 'Uruguay Primera División'
$ws.Range('D236').Value = 'Uruguay Clausura'
$ws.Range('E236').Value = 45267.70833333334
$ws.Range('F236').Value = 'Nacional De Football'
$ws.Range('G236').Value = 'Torque'
$ws.Range('H236').Value = 1
$ws.Range('I236').Value = 1
$ws.Range('J236').Value = 'D'
$ws.Range('K236').Value = 1.666
$ws.Range('L236').Value = 3.9
$ws.Range('M236').Value = 4.5
$ws.Range('N236').Value = 1.615
$ws.Range('O236').Value = 4
$ws.Range('P236').Value = 4.75
$ws.Range('Q236').Value = -0.75
$ws.Range('R236').Value = 1.8
$ws.Range('S236').Value = 2.05
$ws.Range('T236').Value = 2.75
$ws.Range('U236').Value = 1.95
$ws.Range('V236').Value = 1.9
$ws.Range('W236').Value = -1
$ws.Range('X236').Value = 3
$ws.Range('Y236').Value = -1
$ws.Range('Z236').Value = -1
$ws.Range('AA236').Value = 1.05
$ws.Range('AB236').Value = -1
$ws.Range('AC236').Value = 0.8999999999999999

# Row 238
$ws.Range('A238').Value = 236
$ws.Range('B238').Value = 7013886
$ws.Range('C238').Value = 'Uruguay Primera División'
$ws.Range('D238').Value = 'Uruguay Clausura'
$ws.Range('E238').Value = 45267.70833333334
$ws.Range('F238').Value = 'Racing Club de Montevideo'
$ws.Range('G238').Value = 'Cerro'
$ws.Range('H238').Value = 0
$ws.Range('I238').Value = 1
$ws.Range('J238').Value = 'A'
$ws.Range('K238').Value = 2.25
$ws.Range('L238').Value = 3.1
$ws.Range('M238').Value = 3.25
$ws.Range('N238').Value = 2.25
$ws.Range('O238').Value = 2.875
$ws.Range('P238').Value = 3.5
$ws.Range('Q238').Value = -0.25
$ws.Range('R238').Value = 1.95
$ws.Range('S238').Value = 1.9
$ws.Range('T238').Value = 2
$ws.Range('U238').Value = 1.925
$ws.Range('V238').Value = 1.925
$ws.Range('W238').Value = -1
$ws.Range('X238').Value = -1
$ws.Range('Y238').Value = 2.5
$ws.Range('Z238').Value = -1
$ws.Range('AA238').Value = 0.8999999999999999
$ws.Range('AB238').Value = -1
$ws.Range('AC238').Value = 0.925

# Row 239
$ws.Range('A239').Value = 237
$ws.Range('B239').Value = 7013885
$ws.Range('C239').Value = 'Uruguay Primera División'
$ws.Range('D239').Value = 'Uruguay Clausura'
$ws.Range('E239').Value = 45267.70833333334
$ws.Range('F239').Value = 'La Luz'
$ws.Range('G239').Value = 'Atletico Fenix Montevideo'
$ws.Range('H239').Value = 0
$ws.Range('I239').Value = 2
$ws.Range('J239').Value = 'A'
$ws.Range('K239').Value = 3
$ws.Range('L239').Value = 3
$ws.Range('M239').Value = 2.4
$ws.Range('N239').Value = 2.9
$ws.Range('O239').Value = 2.75
$ws.Range('P239').Value = 2.6
$ws.Range('Q239').Value = 0
$ws.Range('R239').Value = 2.025
$ws.Range('S239').Value = 1.825
$ws.Range('T239').Value = 2
$ws.Range('U239').Value = 2.025
$ws.Range('V239').Value = 1.825
$ws.Range('W239').Value = -1
$ws.Range('X239').Value = -1
$ws.Range('Y239').Value = 1.6
$ws.Range('Z239').Value = -1
$ws.Range('AA239').Value = 0.825
$ws.Range('AB239').Value = 0
$ws.Range('AC239').Value = -0

# Row 241
$ws.Range('A241').Value = 239
$ws.Range('B241').Value = 7797510
$ws.Range('C241').Value = 'Uruguay Primera División'
$ws.Range('D241').Value = 'Uruguay Apertura'
$ws.Range('E241').Value = 45338.85416666666
$ws.Range('F241').Value = 'Nacional De Football'
$ws.Range('G241').Value = 'CA River Plate'
$ws.Range('H241').Value = 2
$ws.Range('I241').Value = 1
$ws.Range('J241').Value = 'H'
$ws.Range('K241').Value = 1.6
$ws.Range('L241').Value = 4
$ws.Range('M241').Value = 5.5
$ws.Range('N241').Value = 1.55
$ws.Range('O241').Value = 4
$ws.Range('P241').Value = 6
$ws.Range('Q241').Value = -1
$ws.Range('R241').Value = 2
$ws.Range('S241').Value = 1.85
$ws.Range('T241').Value = 2.25
$ws.Range('U241').Value = 1.825
$ws.Range('V241').Value = 2.025
$ws.Range('W241').Value = 0.55
$ws.Range('X241').Value = -1
$ws.Range('Y241').Value = -1
$ws.Range('Z241').Value = 0
$ws.Range('AA241').Value = -0
$ws.Range('AB241').Value = 0.825
$ws.Range('AC241').Value = -1

# Row 242
$ws.Range('A242').Value = 240
$ws.Range('B242').Value = 7797528
$ws.Range('C242').Value = 'Uruguay Primera División'
$ws.Range('D242').Value = 'Uruguay Apertura'
$ws.Range('E242').Value = 45339.70833333334
$ws.Range('F242').Value = 'Atletico Fenix Montevideo'
$ws.Range('G242').Value = 'Danubio'
$ws.Range('H242').Value = 1
$ws.Range('I242').Value = 2
$ws.Range('J242').Value = 'A'
$ws.Range('K242').Value = 2.6
$ws.Range('L242').Value = 2.9
$ws.Range('M242').Value = 2.9
$ws.Range('N242').Value = 2.75
$ws.Range('O242').Value = 2.875
$ws.Range('P242').Value = 2.75
$ws.Range('Q242').Value = 0
$ws.Range('R242').Value = 1.9
$ws.Range('S242').Value = 1.95
$ws.Range('T242').Value = 2
$ws.Range('U242').Value = 2.1
$ws.Range('V242').Value = 1.775
$ws.Range('W242').Value = -1
$ws.Range('X242').Value = -1
$ws.Range('Y242').Value = 1.75
$ws.Range('Z242').Value = -1
$ws.Range('AA242').Value = 0.95
$ws.Range('AB242').Value = 1.1
$ws.Range('AC242').Value = -1

# Row 243
$ws.Range('A243').Value = 241
$ws.Range('B243').Value = 7797532
$ws.Range('C243').Value = 'Uruguay Primera División'
$ws.Range('D243').Value = 'Uruguay Apertura'
$ws.Range('E243').Value = 45339.80208333334
$ws.Range('F243').Value = 'Miramar Misiones'
$ws.Range('G243').Value = 'Club Atletico Progreso'
$ws.Range('H243').Value = 2
$ws.Range('I243').Value = 3
$ws.Range('J243').Value = 'A'
$ws.Range('K243').Value = 2.5
$ws.Range('L243').Value = 3.2
$ws.Range('M243').Value = 2.8
$ws.Range('N243').Value = 2.4
$ws.Range('O243').Value = 3.2
$ws.Range('P243').Value = 2.9
$ws.Range('Q243').Value = 0
$ws.Range('R243').Value = 1.725
$ws.Range('S243').Value = 2.075
$ws.Range('T243').Value = 2
$ws.Range('U243').Value = 1.825
$ws.Range('V243').Value = 2.025
$ws.Range('W243').Value = -1
$ws.Range('X243').Value = -1
$ws.Range('Y243').Value = 1.9
$ws.Range('Z243').Value = -1
$ws.Range('AA243').Value = 1.075
$ws.Range('AB243').Value = 0.825
$ws.Range('AC243').Value = -1

# Row 244
$ws.Range('A244').Value = 242
$ws.Range('B244').Value = 7797529
$ws.Range('C244').Value = 'Uruguay Primera División'
$ws.Range('D244').Value = 'Uruguay Apertura'
$ws.Range('E244').Value = 45339.89583333334
$ws.Range('F244').Value = 'Deportivo Maldonado'
$ws.Range('G244').Value = 'Boston River'
$ws.Range('H244').Value = 1
$ws.Range('I244').Value = 2
$ws.Range('J244').Value = 'A'
$ws.Range('K244').Value = 2.3
$ws.Range('L244').Value = 3.2
$ws.Range('M244').Value = 3.1
$ws.Range('N244').Value = 2.25
$ws.Range('O244').Value = 3.2
$ws.Range('P244').Value = 3.25
$ws.Range('Q244').Value = -0.25
$ws.Range('R244').Value = 1.95
$ws.Range('S244').Value = 1.9
$ws.Range('T244').Value = 2.25
$ws.Range('U244').Value = 1.95
$ws.Range('V244').Value = 1.9
$ws.Range('W244').Value = -1
$ws.Range('X244').Value = -1
$ws.Range('Y244').Value = 2.25
$ws.Range('Z244').Value = -1
$ws.Range('AA244').Value = 0.8999999999999999
$ws.Range('AB244').Value = 0.95
$ws.Range('AC244').Value = -1

# Row 245
$ws.Range('A245').Value = 243
$ws.Range('B245').Value = 7797530
$ws.Range('C245').Value = 'Uruguay Primera División'
$ws.Range('D245').Value = 'Uruguay Apertura'
$ws.Range('E245').Value = 45340.41666666666
$ws.Range('F245').Value = 'Cerro'
$ws.Range('G245').Value = 'Montevideo Wanderers'
$ws.Range('H245').Value = 1
$ws.Range('I245').Value = 1
$ws.Range('J245').Value = 'D'
$ws.Range('K245').Value = 2.4
$ws.Range('L245').Value = 3
$ws.Range('M245').Value = 3.25
$ws.Range('N245').Value = 2.625
$ws.Range('O245').Value = 3
$ws.Range('P245').Value = 3
$ws.Range('Q245').Value = 0
$ws.Range('R245').Value = 1.775
$ws.Range('S245').Value = 2.1
$ws.Range('T245').Value = 2
$ws.Range('U245').Value = 2
$ws.Range('V245').Value = 1.85
$ws.Range('W245').Value = -1
$ws.Range('X245').Value = 2
$ws.Range('Y245').Value = -1
$ws.Range('Z245').Value = 0
$ws.Range('AA245').Value = -0
$ws.Range('AB245').Value = 0
$ws.Range('AC245').Value = -0

# Row 246
$ws.Range('A246').Value = 244
$ws.Range('B246').Value = 7796575
$ws.Range('C246').Value = 'Uruguay Primera División'
$ws.Range('D246').Value = 'Uruguay Apertura'
$ws.Range('E246').Value = 45340.70833333334
$ws.Range('F246').Value = 'Racing Club de Montevideo'
$ws.Range('G246').Value = 'Liverpool Montevideo'
$ws.Range('H246').Value = 2
$ws.Range('I246').Value = 2
$ws.Range('J246').Value = 'D'
$ws.Range('K246').Value = 3.8
$ws.Range('L246').Value = 3.3
$ws.Range('M246').Value = 1.95
$ws.Range('N246').Value = 3
$ws.Range('O246').Value = 3.2
$ws.Range('P246').Value = 2.4
$ws.Range('Q246').Value = 0.25
$ws.Range('R246').Value = 1.775
$ws.Range('S246').Value = 2.1
$ws.Range('T246').Value = 2.25
$ws.Range('U246').Value = 2.025
$ws.Range('V246').Value = 1.825
$ws.Range('W246').Value = -1
$ws.Range('X246').Value = 2.2
$ws.Range('Y246').Value = -1
$ws.Range('Z246').Value = 0.3875
$ws.Range('AA246').Value = -0.5
$ws.Range('AB246').Value = 1.025
$ws.Range('AC246').Value = -1

# Row 247
$ws.Range('A247').Value = 245
$ws.Range('B247').Value = 7797533
$ws.Range('C247').Value = 'Uruguay Primera División'
$ws.Range('D247').Value = 'Uruguay Apertura'
$ws.Range('E247').Value = 45340.83333333334
$ws.Range('F247').Value = 'Cerro Largo'
$ws.Range('G247').Value = 'Penarol'
$ws.Range('H247').Value = 1
$ws.Range('I247').Value = 2
$ws.Range('J247').Value = 'A'
$ws.Range('K247').Value = 5.25
$ws.Range('L247').Value = 3.5
$ws.Range('M247').Value = 1.7
$ws.Range('N247').Value = 5.75
$ws.Range('O247').Value = 3.4
$ws.Range('P247').Value = 1.666
$ws.Range('Q247').Value = 0.75
$ws.Range('R247').Value = 1.9
$ws.Range('S247').Value = 1.95
$ws.Range('T247').Value = 2
$ws.Range('U247').Value = 1.95
$ws.Range('V247').Value = 1.9
$ws.Range('W247').Value = -1
$ws.Range('X247').Value = -1
$ws.Range('Y247').Value = 0.6659999999999999
$ws.Range('Z247').Value = -0.5
$ws.Range('AA247').Value = 0.475
$ws.Range('AB247').Value = 0.95
$ws.Range('AC247').Value = -1

# --- Create new rows (copy style from row 241 template for A/E columns) ---
# Row 248
$ws.Range('A241').Copy($ws.Range('A248'))
$ws.Range('E241').Copy($ws.Range('E248'))
$ws.Range('A248').Value = 246
$ws.Range('B248').Value = 7825144
$ws.Range('C248').Value = 'Uruguay Primera División'
$ws.Range('D248').Value = 'Uruguay Apertura'
$ws.Range('E248').Value = 45345.70833333334
$ws.Range('F248').Value = 'CA River Plate'
$ws.Range('G248').Value = 'Deportivo Maldonado'
$ws.Range('H248').Value = 3
$ws.Range('I248').Value = 1
$ws.Range('J248').Value = 'H'
$ws.Range('K248').Value = 2.375
$ws.Range('L248').Value = 3.1
$ws.Range('M248').Value = 3
$ws.Range('N248').Value = 2.375
$ws.Range('O248').Value = 3.1
$ws.Range('P248').Value = 3
$ws.Range('Q248').Value = -0.25
$ws.Range('R248').Value = 2.05
$ws.Range('S248').Value = 1.8
$ws.Range('T248').Value = 2.25
$ws.Range('U248').Value = 2.025
$ws.Range('V248').Value = 1.825
$ws.Range('W248').Value = 1.375
$ws.Range('X248').Value = -1
$ws.Range('Y248').Value = -1
$ws.Range('Z248').Value = 1.05
$ws.Range('AA248').Value = -1
$ws.Range('AB248').Value = 1.025
$ws.Range('AC248').Value = -1

# Row 249
$ws.Range('A241').Copy($ws.Range('A249'))
$ws.Range('E241').Copy($ws.Range('E249'))
$ws.Range('A249').Value = 247
$ws.Range('B249').Value = 7825143
$ws.Range('C249').Value = 'Uruguay Primera División'
$ws.Range('D249').Value = 'Uruguay Apertura'
$ws.Range('E249').Value = 45346.70833333334
$ws.Range('F249').Value = 'Boston River'
$ws.Range('G249').Value = 'Danubio'
$ws.Range('K249').Value = 2.75
$ws.Range('L249').Value = 3.1
$ws.Range('M249').Value = 2.5
$ws.Range('N249').Value = 2.55
$ws.Range('O249').Value = 3.1
$ws.Range('P249').Value = 2.7
$ws.Range('Q249').Value = 0
$ws.Range('R249').Value = 1.875
$ws.Range('S249').Value = 1.975
$ws.Range('T249').Value = 2.25
$ws.Range('U249').Value = 1.95
$ws.Range('V249').Value = 1.9
$ws.Range('W249').Value = 0
$ws.Range('X249').Value = 0
$ws.Range('Y249').Value = 0
$ws.Range('Z249').Value = 0
$ws.Range('AA249').Value = 0

# Row 250
$ws.Range('A241').Copy($ws.Range('A250'))
$ws.Range('E241').Copy($ws.Range('E250'))
$ws.Range('A250').Value = 248
$ws.Range('B250').Value = 7825103
$ws.Range('C250').Value = 'Uruguay Primera División'
$ws.Range('D250').Value = 'Uruguay Apertura'
$ws.Range('E250').Value = 45346.83333333334
$ws.Range('F250').Value = 'Penarol'
$ws.Range('G250').Value = 'Miramar Misiones'
$ws.Range('K250').Value = 1.4
$ws.Range('L250').Value = 4
$ws.Range('M250').Value = 9
$ws.Range('N250').Value = 1.363
$ws.Range('O250').Value = 4.333
$ws.Range('P250').Value = 9
$ws.Range('Q250').Value = -1.25
$ws.Range('R250').Value = 1.875
$ws.Range('S250').Value = 1.975
$ws.Range('T250').Value = 2.5
$ws.Range('U250').Value = 1.975
$ws.Range('V250').Value = 1.875
$ws.Range('W250').Value = 0
$ws.Range('X250').Value = 0
$ws.Range('Y250').Value = 0
$ws.Range('Z250').Value = 0
$ws.Range('AA250').Value = 0

# Row 251
$ws.Range('A241').Copy($ws.Range('A251'))
$ws.Range('E241').Copy($ws.Range('E251'))
$ws.Range('A251').Value = 249
$ws.Range('B251').Value = 7825146
$ws.Range('C251').Value = 'Uruguay Primera División'
$ws.Range('D251').Value = 'Uruguay Apertura'
$ws.Range('E251').Value = 45347.41666666666
$ws.Range('F251').Value = 'Cerro Largo'
$ws.Range('G251').Value = 'Atletico Fenix Montevideo'
$ws.Range('K251').Value = 2.625
$ws.Range('L251').Value = 2.875
$ws.Range('M251').Value = 2.75
$ws.Range('N251').Value = 2.45
$ws.Range('O251').Value = 2.875
$ws.Range('P251').Value = 3
$ws.Range('Q251').Value = -0.25
$ws.Range('R251').Value = 2.125
$ws.Range('S251').Value = 1.75
$ws.Range('T251').Value = 2
$ws.Range('U251').Value = 1.975
$ws.Range('V251').Value = 1.875
$ws.Range('W251').Value = 0
$ws.Range('X251').Value = 0
$ws.Range('Y251').Value = 0
$ws.Range('Z251').Value = 0
$ws.Range('AA251').Value = 0

# Row 252
$ws.Range('A241').Copy($ws.Range('A252'))
$ws.Range('E241').Copy($ws.Range('E252'))
$ws.Range('A252').Value = 250
$ws.Range('B252').Value = 7825104
$ws.Range('C252').Value = 'Uruguay Primera División'
$ws.Range('D252').Value = 'Uruguay Apertura'
$ws.Range('E252').Value = 45347.70833333334
$ws.Range('F252').Value = 'Liverpool Montevideo'
$ws.Range('G252').Value = 'Nacional De Football'
$ws.Range('K252').Value = 2.875
$ws.Range('L252').Value = 3.3
$ws.Range('M252').Value = 2.375
$ws.Range('N252').Value = 2.75
$ws.Range('O252').Value = 3.3
$ws.Range('P252').Value = 2.45
$ws.Range('Q252').Value = 0
$ws.Range('R252').Value = 2.05
$ws.Range('S252').Value = 1.8
$ws.Range('T252').Value = 2.5
$ws.Range('U252').Value = 1.925
$ws.Range('V252').Value = 1.925
$ws.Range('W252').Value = 0
$ws.Range('X252').Value = 0
$ws.Range('Y252').Value = 0
$ws.Range('Z252').Value = 0
$ws.Range('AA252').Value = 0

# Row 253
$ws.Range('A241').Copy($ws.Range('A253'))
$ws.Range('E241').Copy($ws.Range('E253'))
$ws.Range('A253').Value = 251
$ws.Range('B253').Value = 7825147
$ws.Range('C253').Value = 'Uruguay Primera División'
$ws.Range('D253').Value = 'Uruguay Apertura'
$ws.Range('E253').Value = 45347.8125
$ws.Range('F253').Value = 'Defensor Sporting'
$ws.Range('G253').Value = 'Cerro'
$ws.Range('K253').Value = 1.727
$ws.Range('L253').Value = 3.6
$ws.Range('M253').Value = 4.333
$ws.Range('N253').Value = 1.7
$ws.Range('O253').Value = 3.6
$ws.Range('P253').Value = 4.5
$ws.Range('Q253').Value = -0.75
$ws.Range('R253').Value = 1.975
$ws.Range('S253').Value = 1.875
$ws.Range('T253').Value = 2.5
$ws.Range('U253').Value = 2.05
$ws.Range('V253').Value = 1.8
$ws.Range('W253').Value = 0
$ws.Range('X253').Value = 0
$ws.Range('Y253').Value = 0
$ws.Range('Z253').Value = 0
$ws.Range('AA253').Value = 0
